$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Week 4" column header to include its date range, like the
# other week columns already have.
$ws.Range("E1").Value = "Week 4 -- May 23 - 29"

# Fill in the Week 4 actuals for each athlete (previously all placeholder 0s).
$ws.Range("E2").Value = 112.9
$ws.Range("E3").Value = 253.6
$ws.Range("E4").Value = 48.4
$ws.Range("E5").Value = 178.8
$ws.Range("E6").Value = 288.8
$ws.Range("E7").Value = 284.4
$ws.Range("E8").Value = 126.3
$ws.Range("E9").Value = 194.2
$ws.Range("E10").Value = 25
$ws.Range("E11").Value = 126.9

# E2 loses its inherited numeric style, matching the rest of row 2 which
# already uses the default "Normal" style for columns C/D.
$ws.Range("E2").Style = "Normal"

# Widen column E to match column D now that its header text is longer.
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(4).ColumnWidth

# Leave the selection on E2 (also resets the scrolled view back to the
# top-left of the sheet).
$ws.Range("E2").Select() | Out-Null
